$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data for observations 21 and 22 got re-sorted upstream, so the
# rows effectively swap contents (the row numbers / positions stay put).
# Only the columns that actually differ between the two rows need touching:
# Id (A), Taxonsorteringsordning (B), TaxonId (E), Artnamn (F),
# Vetenskapligt namn (G), Auktor (H), Ost (Q), Nord (R), Starttid (Z),
# Sluttid (AB).

$cols = @(1, 2, 5, 6, 7, 8, 17, 18, 26, 28)

foreach ($c in $cols) {
    $cell21 = $ws.Cells.Item(21, $c)
    $cell22 = $ws.Cells.Item(22, $c)

    $v21 = $cell21.Value2
    $v22 = $cell22.Value2

    $cell21.Value2 = $v22
    $cell22.Value2 = $v21
}
